$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "'2026-02-02"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = '산업'
$ws.Range("C2").Value = '대한민국 AI 풀스택, 사우디 시장 진출 본격화'
$ws.Range("D2").Value = "'2026-02-02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = 'https://n.news.naver.com/mnews/article/031/0001001993?sid=105'

# Row 3
$ws.Range("A3").Value = "'2026-02-02"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = '산업'
$ws.Range("C3").Value = '주가조작 초동 대응 강화…거래소, AI 시장감시 체계 가동'
$ws.Range("D3").Value = "'2026-02-02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = 'https://n.news.naver.com/mnews/article/003/0013743257?sid=101'

# Row 4
$ws.Range("A4").Value = "'2026-02-02"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = '기업'
$ws.Range("C4").Value = '오픈AI ''투자 보류설'' 정면 반박한 젠슨 황'
$ws.Range("D4").Value = "'2026-02-01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = 'https://n.news.naver.com/mnews/article/015/0005245153?sid=105'

# Row 5
$ws.Range("A5").Value = "'2026-02-02"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = '기술'
$ws.Range("C5").Value = '금감원, 가상자산 시세조종 AI로 적발…자동탐지 알고리즘 도입'
$ws.Range("D5").Value = "'2026-02-02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = 'https://n.news.naver.com/mnews/article/001/0015879611?sid=101'

# Row 6
$ws.Range("A6").Value = "'2026-02-02"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = '정책'
$ws.Range("C6").Value = '산업부, 중견기업 R&D에 655억 투입…"지역 발전·AI 혁신 지원"'
$ws.Range("D6").Value = "'2026-02-02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = 'https://n.news.naver.com/mnews/article/001/0015879389?sid=101'

# Row 7
$ws.Range("A7").Value = "'2026-02-02"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = '기술'
$ws.Range("C7").Value = '토스증권, AI 어닝콜 PC버전 출시'
$ws.Range("D7").Value = "'2026-02-02"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = 'https://sports.hankooki.com/news/articleView.html?idxno=6923817'

# Row 8
$ws.Range("A8").Value = "'2026-02-02"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = '기업'
$ws.Range("C8").Value = '젠슨 황, 오픈AI ''불만설''에 입 열어…"대규모 투자 진행"'
$ws.Range("D8").Value = "'2026-02-01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = 'https://n.news.naver.com/mnews/article/015/0005245186?sid=101'

# Row 9
$ws.Range("A9").Value = "'2026-02-02"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = '정책'
$ws.Range("C9").Value = '과기부 "한국 AI 풀스택, 사우디 진출 본격화"'
$ws.Range("D9").Value = "'2026-02-02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = 'https://n.news.naver.com/mnews/article/215/0001240176?sid=101'

# Row 10
$ws.Range("A10").Value = "'2026-02-02"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = '정부(과기부)'
$ws.Range("C10").Value = '가비아, ‘AX 지원 프로모션’…AI 서비스 크레딧 지원'
$ws.Range("D10").Value = "'2026-02-02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = 'http://www.datanews.co.kr/news/article.html?no=143146'

# Row 11
$ws.Range("A11").Value = "'2026-02-02"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = '정부(과기부)'
$ws.Range("C11").Value = '한양대 ERICA AI융합연구소, ''연구개발 및 인재 양성 강화 업무협약'' 체...'
$ws.Range("D11").Value = "'2026-02-02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = 'https://www.kfenews.co.kr/news/articleView.html?idxno=653866'

# Row 12
$ws.Range("A12").Value = "'2026-02-02"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = '정부(과기부)'
$ws.Range("C12").Value = '대덕특구 연구소기업 ''한다랩'', 나스닥 첫 진출 성공'
$ws.Range("D12").Value = "'2026-02-02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = 'http://www.enewstoday.co.kr/news/articleView.html?idxno=2389438'

# Row 13
$ws.Range("A13").Value = "'2026-02-02"
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Value = '정부(과기부)'
$ws.Range("C13").Value = '한국정보인증, IITP ''정보보호 핵심 원천기술 개발사업'' 4차년도 과제 성...'
$ws.Range("D13").Value = "'2026-02-02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = 'https://n.news.naver.com/mnews/article/022/0004102540?sid=101'

# Row 14
$ws.Range("A14").Value = "'2026-02-02"
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").Value = '정부(과기부)'
$ws.Range("C14").Value = '과기부 "한국 AI 풀스택, 사우디 진출 본격화"'
$ws.Range("D14").Value = "'2026-02-02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = 'https://n.news.naver.com/mnews/article/215/0001240176?sid=101'
